# Add 2022-Q4 data
#
# Before: 总计, 2022-Q3, 2022-Q1
# After : 总计, 2022-Q4, 2022-Q3, 2022-Q1
#
# The existing "2022-Q3" sheet is duplicated (so the brand-new "2022-Q3"
# sheet keeps the untouched Q3 figures/formatting), the original sheet is
# renamed to "2022-Q4" and refreshed with the new quarter's numbers, and the
# "总计" summary sheet gets a new leading row for 2022-Q4 with everything
# else shifted down one slot.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate the current "2022-Q3" sheet right after itself. The copy
#    inherits all values/styles, so it becomes the new, permanent
#    "2022-Q3" sheet (holding the untouched Q3 figures).
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($null, $wsQ3)
$wsQ3Dup = $wb.Worksheets.Item("2022-Q3 (2)")

# Free up the "2022-Q3" name on the original sheet before renaming, then
# give the duplicate the permanent "2022-Q3" name.
$wsQ3.Name = "2022-Q4"
$wsQ3Dup.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2. Refresh the figures on the (renamed) "2022-Q4" sheet. These columns
#    are stored as text in the workbook, so force text by entering the
#    values with a leading apostrophe, then reset the cell style back to
#    Normal (the apostrophe trick otherwise marks the cell "quote
#    prefixed", which nudges the style index).
# ---------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Item("2022-Q4")

$wsQ4.Range("D2").Value = "'27.13"
$wsQ4.Range("D2").Style = "Normal"

$wsQ4.Range("E2").Value = "'99.24"
$wsQ4.Range("E2").Style = "Normal"

$wsQ4.Range("F2").Value = "'2.80"
$wsQ4.Range("F2").Style = "Normal"

$wsQ4.Range("G2").Value = "'0.7596"
$wsQ4.Range("G2").Style = "Normal"

# ---------------------------------------------------------------------
# 3. Rewrite the "总计" summary sheet: shift the Q3/Q1 rows down one row
#    and insert the new Q4 row on top.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A4").PasteSpecial(-4122)
$wsTotal.Range("B4").Value = "2022-Q1"
$wsTotal.Range("C4").Value = 1
$wsTotal.Range("D4").Value = 0.62

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.75

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.76

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. The sheet-copy above leaves the duplicate as the active tab; restore
#    the original active sheet ("总计" / the workbook's first tab).
# ---------------------------------------------------------------------
$wsTotal.Activate()
